# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") holds strikeout counts for each outing row (rows 2-12).
# The save_data regeneration recomputed these values (K instead of Strike#).
$kValues = @{
    2  = 0
    3  = 2
    4  = 3
    5  = 3
    6  = 0
    7  = 5
    8  = 4
    9  = 1
    10 = 0
    11 = 1
    12 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
